$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 57.38695133333332
$ws.Range("H2").Value = 172.160854
$ws.Range("I2").Value = 0.6542464432660997
$ws.Range("J2").Value = 0.6542464432660998
$ws.Range("M2").Value = 7.368211
$ws.Range("N2").Value = 22.104633
$ws.Range("O2").Value = 0.1280150446959744
$ws.Range("P2").Value = 0.1280150446959744
$ws.Range("Q2").Value = 422.8391660707312
$ws.Range("R2").Value = 3805.552494636581
$ws.Range("S2").Value = 0.08375338767689205
$ws.Range("T2").Value = 0.08375338767689207
$ws.Range("G3").Value = 57.38695133333332
$ws.Range("H3").Value = 172.160854
$ws.Range("I3").Value = 0.6542464432660997
$ws.Range("J3").Value = 0.6542464432660998
$ws.Range("M3").Value = 18.05628333333334
$ws.Range("N3").Value = 54.16885000000001
$ws.Range("O3").Value = 0.3137092461059877
$ws.Range("P3").Value = 0.3137092461059876
$ws.Range("Q3").Value = 1036.195052910878
$ws.Range("R3").Value = 9325.755476197899
$ws.Range("S3").Value = 0.205243158484532
$ws.Range("T3").Value = 0.205243158484532
$ws.Range("G4").Value = 57.38695133333332
$ws.Range("H4").Value = 172.160854
$ws.Range("I4").Value = 0.6542464432660997
$ws.Range("J4").Value = 0.6542464432660998
$ws.Range("M4").Value = 18.794071
$ws.Range("N4").Value = 56.38221299999999
$ws.Range("O4").Value = 0.3265275436716344
$ws.Range("P4").Value = 0.3265275436716343
$ws.Range("Q4").Value = 1078.534437832211
$ws.Range("R4").Value = 9706.809940489899
$ws.Range("S4").Value = 0.2136294840755829
$ws.Range("T4").Value = 0.2136294840755829
$ws.Range("G5").Value = 57.38695133333332
$ws.Range("H5").Value = 172.160854
$ws.Range("I5").Value = 0.6542464432660997
$ws.Range("J5").Value = 0.6542464432660998
$ws.Range("M5").Value = 13.338818
$ws.Range("N5").Value = 40.016454
$ws.Range("O5").Value = 0.2317481655264036
$ws.Range("P5").Value = 0.2317481655264036
$ws.Range("Q5").Value = 765.4740994101904
$ws.Range("R5").Value = 6889.266894691714
$ws.Range("S5").Value = 0.1516204130290929
$ws.Range("T5").Value = 0.151620413029093
$ws.Range("I6").Value = 0.16357689713892
$ws.Range("J6").Value = 0.16357689713892
$ws.Range("M6").Value = 7.368211
$ws.Range("N6").Value = 22.104633
$ws.Range("O6").Value = 0.1280150446959744
$ws.Range("P6").Value = 0.1280150446959744
$ws.Range("Q6").Value = 105.7196710605987
$ws.Range("R6").Value = 951.4770395453882
$ws.Range("S6").Value = 0.02094030379846765
$ws.Range("T6").Value = 0.02094030379846765
$ws.Range("I7").Value = 0.16357689713892
$ws.Range("J7").Value = 0.16357689713892
$ws.Range("M7").Value = 18.05628333333334
$ws.Range("N7").Value = 54.16885000000001
$ws.Range("O7").Value = 0.3137092461059877
$ws.Range("P7").Value = 0.3137092461059876
$ws.Range("Q7").Value = 259.072973694289
$ws.Range("S7").Value = 0.05131558508180727
$ws.Range("T7").Value = 0.05131558508180726
$ws.Range("I8").Value = 0.16357689713892
$ws.Range("J8").Value = 0.16357689713892
$ws.Range("M8").Value = 18.794071
$ws.Range("N8").Value = 56.38221299999999
$ws.Range("O8").Value = 0.3265275436716344
$ws.Range("P8").Value = 0.3265275436716343
$ws.Range("Q8").Value = 269.6588091749186
$ws.Range("R8").Value = 2426.929282574268
$ws.Range("S8").Value = 0.05341236242419913
$ws.Range("T8").Value = 0.05341236242419913
$ws.Range("I9").Value = 0.16357689713892
$ws.Range("J9").Value = 0.16357689713892
$ws.Range("M9").Value = 13.338818
$ws.Range("N9").Value = 40.016454
$ws.Range("O9").Value = 0.2317481655264036
$ws.Range("P9").Value = 0.2317481655264036
$ws.Range("Q9").Value = 191.3864099843493
$ws.Range("R9").Value = 1722.477689859144
$ws.Range("S9").Value = 0.03790864583444593
$ws.Range("T9").Value = 0.03790864583444593
$ws.Range("G10").Value = 14.516389
$ws.Range("H10").Value = 43.549167
$ws.Range("I10").Value = 0.1654957381714162
$ws.Range("J10").Value = 0.1654957381714162
$ws.Range("M10").Value = 7.368211
$ws.Range("N10").Value = 22.104633
$ws.Range("O10").Value = 0.1280150446959744
$ws.Range("P10").Value = 0.1280150446959744
$ws.Range("Q10").Value = 106.959817110079
$ws.Range("R10").Value = 962.6383539907109
$ws.Range("S10").Value = 0.02118594431900712
$ws.Range("T10").Value = 0.02118594431900712
$ws.Range("G11").Value = 14.516389
$ws.Range("H11").Value = 43.549167
$ws.Range("I11").Value = 0.1654957381714162
$ws.Range("J11").Value = 0.1654957381714162
$ws.Range("M11").Value = 18.05628333333334
$ws.Range("N11").Value = 54.16885000000001
$ws.Range("O11").Value = 0.3137092461059877
$ws.Range("P11").Value = 0.3137092461059876
$ws.Range("Q11").Value = 262.1120327608834
$ws.Range("R11").Value = 2359.00829484795
$ws.Range("S11").Value = 0.05191754325550889
$ws.Range("T11").Value = 0.05191754325550888
$ws.Range("G12").Value = 14.516389
$ws.Range("H12").Value = 43.549167
$ws.Range("I12").Value = 0.1654957381714162
$ws.Range("J12").Value = 0.1654957381714162
$ws.Range("M12").Value = 18.794071
$ws.Range("N12").Value = 56.38221299999999
$ws.Range("O12").Value = 0.3265275436716344
$ws.Range("P12").Value = 0.3265275436716343
$ws.Range("Q12").Value = 272.8220455296189
$ws.Range("R12").Value = 2455.398409766571
$ws.Range("S12").Value = 0.05403891687323646
$ws.Range("T12").Value = 0.05403891687323645
$ws.Range("G13").Value = 14.516389
$ws.Range("H13").Value = 43.549167
$ws.Range("I13").Value = 0.1654957381714162
$ws.Range("J13").Value = 0.1654957381714162
$ws.Range("M13").Value = 13.338818
$ws.Range("N13").Value = 40.016454
$ws.Range("O13").Value = 0.2317481655264036
$ws.Range("P13").Value = 0.2317481655264036
$ws.Range("Q13").Value = 193.631470888202
$ws.Range("R13").Value = 1742.683237993818
$ws.Range("S13").Value = 0.03835333372366371
$ws.Range("T13").Value = 0.03835333372366371
$ws.Range("G14").Value = 1.46316
$ws.Range("H14").Value = 4.389480000000001
$ws.Range("I14").Value = 0.01668092142356404
$ws.Range("J14").Value = 0.01668092142356404
$ws.Range("M14").Value = 7.368211
$ws.Range("N14").Value = 22.104633
$ws.Range("O14").Value = 0.1280150446959744
$ws.Range("P14").Value = 0.1280150446959744
$ws.Range("Q14").Value = 10.78087160676
$ws.Range("R14").Value = 97.02784446084002
$ws.Range("S14").Value = 0.002135408901607588
$ws.Range("T14").Value = 0.002135408901607588
$ws.Range("G15").Value = 1.46316
$ws.Range("H15").Value = 4.389480000000001
$ws.Range("I15").Value = 0.01668092142356404
$ws.Range("J15").Value = 0.01668092142356404
$ws.Range("M15").Value = 18.05628333333334
$ws.Range("N15").Value = 54.16885000000001
$ws.Range("O15").Value = 0.3137092461059877
$ws.Range("P15").Value = 0.3137092461059876
$ws.Range("Q15").Value = 26.41923152200001
$ws.Range("R15").Value = 237.7730836980001
$ws.Range("S15").Value = 0.005232959284139492
$ws.Range("T15").Value = 0.005232959284139492
$ws.Range("G16").Value = 1.46316
$ws.Range("H16").Value = 4.389480000000001
$ws.Range("I16").Value = 0.01668092142356404
$ws.Range("J16").Value = 0.01668092142356404
$ws.Range("M16").Value = 18.794071
$ws.Range("N16").Value = 56.38221299999999
$ws.Range("O16").Value = 0.3265275436716344
$ws.Range("P16").Value = 0.3265275436716343
$ws.Range("Q16").Value = 27.49873292436
$ws.Range("R16").Value = 247.48859631924
$ws.Range("S16").Value = 0.005446780298615907
$ws.Range("T16").Value = 0.005446780298615908
$ws.Range("G17").Value = 1.46316
$ws.Range("H17").Value = 4.389480000000001
$ws.Range("I17").Value = 0.01668092142356404
$ws.Range("J17").Value = 0.01668092142356404
$ws.Range("M17").Value = 13.338818
$ws.Range("N17").Value = 40.016454
$ws.Range("O17").Value = 0.2317481655264036
$ws.Range("P17").Value = 0.2317481655264036
$ws.Range("Q17").Value = 19.51682494488
$ws.Range("R17").Value = 175.65142450392
$ws.Range("S17").Value = 0.003865772939201051
$ws.Range("T17").Value = 0.003865772939201052
